$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the previous last data row (70) down onto the new
# row 71 first, so the new cells reuse the existing date/time cell styles
# instead of creating new ones.
$ws.Range("A70:D70").Copy()
$ws.Range("A71:D71").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Add new row 71 with date/time/files/disk_space data.
$ws.Range("A71").Value = 45021
$ws.Range("B71").Value = 0.33241898148148147
$ws.Range("C71").Value = 75499
$ws.Range("D71").Value = 1430

# Update the selected cell to reflect the new last empty row, like Excel
# does automatically after entering data in the previous last row.
$ws.Range("A72").Select()
